$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "B2" "Bitcoin"
Set-TextValue "C2" "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
Set-TextValue "D2" "26.536.56"
Set-TextValue "E2" "  +0.08%  "

Set-TextValue "B3" "Ethereum"
Set-TextValue "C3" "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
Set-TextValue "D3" "1.813.03"
Set-TextValue "E3" "  +0.06%  "

Set-TextValue "B4" "TetherUSD"
Set-TextValue "C4" "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  -0.46%  "

Set-TextValue "B5" "USDC"
Set-TextValue "C5" "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D5" "1.003"
Set-TextValue "E5" "  -0.37%  "

Set-TextValue "B6" "BNB"
Set-TextValue "C6" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue "D6" "305.52"
Set-TextValue "E6" "  -0.89%  "

Set-TextValue "B7" "XRP"
Set-TextValue "C7" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue "D7" "0.4544"
Set-TextValue "E7" "  -0.29%  "

Set-TextValue "B8" "Cardano"
Set-TextValue "C8" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D8" "0.3593"
Set-TextValue "E8" "  -1.95%  "

Set-TextValue "B9" "Dogecoin"
Set-TextValue "C9" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D9" "0.07108"
Set-TextValue "E9" "  -0.35%  "

Set-TextValue "B10" "Polygon"
Set-TextValue "C10" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D10" "0.8925"
Set-TextValue "E10" "  +1.49%  "

Set-TextValue "B11" "TRON"
Set-TextValue "C11" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D11" "0.07728"
Set-TextValue "E11" "  -0.40%  "

Set-TextValue "B12" "Solana"
Set-TextValue "C12" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D12" "19.30"
Set-TextValue "E12" "  -0.33%  "

Set-TextValue "B13" "WrappedEther"
Set-TextValue "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.802.68"
Set-TextValue "E13" "  -0.47%  "

Set-TextValue "B14" "Polkadot"
Set-TextValue "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "5.257"
Set-TextValue "E14" "  -0.61%  "

Set-TextValue "B15" "Chainlink"
Set-TextValue "C15" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D15" "6.301"
Set-TextValue "E15" "  -1.09%  "

Set-TextValue "B16" "Litecoin"
Set-TextValue "C16" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D16" "85.88"
Set-TextValue "E16" "  -0.88%  "

Set-TextValue "B17" "BinanceUSD"
Set-TextValue "C17" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D17" "1.005"
Set-TextValue "E17" "  -0.36%  "

Set-TextValue "B18" "ShibaInu"
Set-TextValue "C18" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D18" "0.000008552"
Set-TextValue "E18" "  -0.42%  "

Set-TextValue "B19" "Dai"
Set-TextValue "C19" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D19" "1.003"
Set-TextValue "E19" "  -0.31%  "

Set-TextValue "B20" "WrappedBTC"
Set-TextValue "C20" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D20" "26.568.53"
Set-TextValue "E20" "  -0.08%  "

Set-TextValue "B21" "Avalanche"
Set-TextValue "C21" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D21" "14.13"
Set-TextValue "E21" "  -0.81%  "

Set-TextValue "B22" "Uniswap"
Set-TextValue "C22" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D22" "4.955"
Set-TextValue "E22" "  -1.08%  "

Set-TextValue "B23" "Cosmos"
Set-TextValue "C23" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D23" "10.50"
Set-TextValue "E23" "  +0.26%  "

Set-TextValue "B24" "Toncoin"
Set-TextValue "C24" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D24" "1.923"
Set-TextValue "E24" "  -3.08%  "

Set-TextValue "B25" "Monero"
Set-TextValue "C25" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D25" "152.06"
Set-TextValue "E25" "  +0.34%  "

Set-TextValue "B26" "EthereumClassic"
Set-TextValue "C26" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D26" "17.79"
Set-TextValue "E26" "  -0.89%  "

Set-TextValue "B27" "LidoDAOToken"
Set-TextValue "C27" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D27" "2.023"
Set-TextValue "E27" "  -2.00%  "

Set-TextValue "B28" "BitcoinCash"
Set-TextValue "C28" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D28" "112.40"
Set-TextValue "E28" "  -0.37%  "

Set-TextValue "B29" "InternetComputer(DFINITY)"
Set-TextValue "C29" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D29" "4.821"
Set-TextValue "E29" "  -0.43%  "

Set-TextValue "B30" "Stellar"
Set-TextValue "C30" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D30" "0.08705"
Set-TextValue "E30" "  +0.28%  "

Set-TextValue "B31" "HuobiToken"
Set-TextValue "C31" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D31" "3.134"
Set-TextValue "E31" "  +2.35%  "

Set-TextValue "B32" "ImmutableX"
Set-TextValue "C32" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D32" "0.7405"
Set-TextValue "E32" "  +1.15%  "

Set-TextValue "B33" "Filecoin"
Set-TextValue "C33" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D33" "4.423"
Set-TextValue "E33" "  -2.21%  "

Set-TextValue "B34" "RenderToken"
Set-TextValue "C34" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D34" "2.710"
Set-TextValue "E34" "  +1.63%  "

Set-TextValue "B35" "ARBITRUM"
Set-TextValue "C35" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D35" "1.109"
Set-TextValue "E35" "  -0.81%  "

Set-TextValue "B36" "TrustWalletToken"
Set-TextValue "C36" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D36" "1.071"
Set-TextValue "E36" "  -1.08%  "

Set-TextValue "B37" "VeChain"
Set-TextValue "C37" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D37" "0.01934"
Set-TextValue "E37" "  -0.94%  "

Set-TextValue "B38" "MXToken"
Set-TextValue "C38" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D38" "2.913"
Set-TextValue "E38" "  +0.53%  "

Set-TextValue "B39" "Hedera"
Set-TextValue "C39" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D39" "0.05081"
Set-TextValue "E39" "  -0.59%  "

Set-TextValue "B40" "TheSandbox"
Set-TextValue "C40" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D40" "0.5086"
Set-TextValue "E40" "  +1.77%  "

Set-TextValue "B41" "FraxShare"
Set-TextValue "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D41" "6.807"
Set-TextValue "E41" "  -2.48%  "

Set-TextValue "B42" "Algorand"
Set-TextValue "C42" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D42" "0.1506"
Set-TextValue "E42" "  -3.49%  "

Set-TextValue "B43" "Aptos"
Set-TextValue "C43" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D43" "8.020"
Set-TextValue "E43" "  -1.90%  "

Set-TextValue "B44" "Decentraland"
Set-TextValue "C44" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D44" "0.4682"
Set-TextValue "E44" "  +1.67%  "

Set-TextValue "B45" "PaxDollar"
Set-TextValue "C45" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D45" "1.003"
Set-TextValue "E45" "  -0.45%  "

Set-TextValue "B46" "EnergySwap"
Set-TextValue "C46" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D46" "9.942"
Set-TextValue "E46" "  -0.60%  "

Set-TextValue "B47" "Quant"
Set-TextValue "C47" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D47" "99.04"
Set-TextValue "E47" "  -1.89%  "

Set-TextValue "B48" "NEARProtocol"
Set-TextValue "C48" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D48" "1.563"
Set-TextValue "E48" "  -1.80%  "

Set-TextValue "B49" "Cronos"
Set-TextValue "C49" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D49" "0.05993"
Set-TextValue "E49" "  -0.09%  "

Set-TextValue "B50" "Aave"
Set-TextValue "C50" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D50" "63.78"
Set-TextValue "E50" "  -0.97%  "

Set-TextValue "B51" "Elrond"
Set-TextValue "C51" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D51" "35.82"
Set-TextValue "E51" "  -1.01%  "

Write-Output "Done"